{"js": "// Replace the date line and the 25 division problems with their updated values.\n// Every source string is unique within the document, so a single\n// search+replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-10-18 Friday\", \"2024-10-19 Saturday\"],\n  [\"33\u00f79=\", \"64\u00f79=\"],\n  [\"63\u00f77=\", \"65\u00f79=\"],\n  [\"55\u00f76=\", \"56\u00f75=\"],\n  [\"80\u00f76=\", \"76\u00f78=\"],\n  [\"98\u00f75=\", \"39\u00f72=\"],\n  [\"19\u00f75=\", \"24\u00f73=\"],\n  [\"35\u00f74=\", \"81\u00f77=\"],\n  [\"66\u00f79=\", \"58\u00f77=\"],\n  [\"85\u00f79=\", \"46\u00f75=\"],\n  [\"66\u00f77=\", \"87\u00f78=\"],\n  [\"40\u00f76=\", \"69\u00f75=\"],\n  [\"52\u00f75=\", \"75\u00f78=\"],\n  [\"14\u00f78=\", \"53\u00f78=\"],\n  [\"87\u00f76=\", \"73\u00f77=\"],\n  [\"12\u00f73=\", \"93\u00f76=\"],\n  [\"68\u00f77=\", \"59\u00f75=\"],\n  [\"52\u00f74=\", \"99\u00f78=\"],\n  [\"23\u00f76=\", \"90\u00f75=\"],\n  [\"20\u00f74=\", \"86\u00f75=\"],\n  [\"24\u00f75=\", \"29\u00f74=\"],\n  [\"63\u00f78=\", \"92\u00f74=\"],\n  [\"49\u00f78=\", \"33\u00f79=\"],\n  [\"56\u00f73=\", \"46\u00f73=\"],\n  [\"53\u00f74=\", \"47\u00f79=\"],\n  [\"54\u00f73=\", \"21\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 division problems with their updated\n# values. Every source string is unique within the document, so a single\n# Find/Replace (replace-all) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-10-18 Friday\", \"2024-10-19 Saturday\"),\n  @(\"33\u00f79=\", \"64\u00f79=\"),\n  @(\"63\u00f77=\", \"65\u00f79=\"),\n  @(\"55\u00f76=\", \"56\u00f75=\"),\n  @(\"80\u00f76=\", \"76\u00f78=\"),\n  @(\"98\u00f75=\", \"39\u00f72=\"),\n  @(\"19\u00f75=\", \"24\u00f73=\"),\n  @(\"35\u00f74=\", \"81\u00f77=\"),\n  @(\"66\u00f79=\", \"58\u00f77=\"),\n  @(\"85\u00f79=\", \"46\u00f75=\"),\n  @(\"66\u00f77=\", \"87\u00f78=\"),\n  @(\"40\u00f76=\", \"69\u00f75=\"),\n  @(\"52\u00f75=\", \"75\u00f78=\"),\n  @(\"14\u00f78=\", \"53\u00f78=\"),\n  @(\"87\u00f76=\", \"73\u00f77=\"),\n  @(\"12\u00f73=\", \"93\u00f76=\"),\n  @(\"68\u00f77=\", \"59\u00f75=\"),\n  @(\"52\u00f74=\", \"99\u00f78=\"),\n  @(\"23\u00f76=\", \"90\u00f75=\"),\n  @(\"20\u00f74=\", \"86\u00f75=\"),\n  @(\"24\u00f75=\", \"29\u00f74=\"),\n  @(\"63\u00f78=\", \"92\u00f74=\"),\n  @(\"49\u00f78=\", \"33\u00f79=\"),\n  @(\"56\u00f73=\", \"46\u00f73=\"),\n  @(\"53\u00f74=\", \"47\u00f79=\"),\n  @(\"54\u00f73=\", \"21\u00f72=\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n\nWrite-Output \"done\"\n"}
